# Auto-generated Excel COM-interop script
# Applies the cryptos.xlsx price/volume update described by the commit diff.
# Each cell keeps its original style/format (text) while only the displayed
# text content changes, matching the source OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("D2")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "53.887.97"
$rng.Style = $origStyle

$rng = $ws.Range("E2")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -4.65%  "
$rng.Style = $origStyle

$rng = $ws.Range("D3")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "2.241.29"
$rng.Style = $origStyle

$rng = $ws.Range("E3")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -5.78%  "
$rng.Style = $origStyle

$rng = $ws.Range("D4")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "0.996"
$rng.Style = $origStyle

$rng = $ws.Range("E4")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -0.63%  "
$rng.Style = $origStyle

$rng = $ws.Range("D5")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "486.91"
$rng.Style = $origStyle

$rng = $ws.Range("E5")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -3.83%  "
$rng.Style = $origStyle

$rng = $ws.Range("D6")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "125.27"
$rng.Style = $origStyle

$rng = $ws.Range("D7")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "0.996"
$rng.Style = $origStyle

$rng = $ws.Range("E7")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -0.28%  "
$rng.Style = $origStyle

$rng = $ws.Range("D8")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "0.523"
$rng.Style = $origStyle

$rng = $ws.Range("E8")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -4.24%  "
$rng.Style = $origStyle

$rng = $ws.Range("D9")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "2.239.13"
$rng.Style = $origStyle

$rng = $ws.Range("E9")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -6.32%  "
$rng.Style = $origStyle

$rng = $ws.Range("D10")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "0.0921"
$rng.Style = $origStyle

$rng = $ws.Range("E10")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -6.69%  "
$rng.Style = $origStyle

$rng = $ws.Range("E11")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -0.17%  "
$rng.Style = $origStyle

$rng = $ws.Range("E12")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -3.16%  "
$rng.Style = $origStyle

$rng = $ws.Range("E13")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -4.88%  "
$rng.Style = $origStyle

$rng = $ws.Range("D14")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "2.611.06"
$rng.Style = $origStyle

$rng = $ws.Range("E14")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -6.77%  "
$rng.Style = $origStyle

$rng = $ws.Range("D15")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "21.23"
$rng.Style = $origStyle

$rng = $ws.Range("E15")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -1.96%  "
$rng.Style = $origStyle

$rng = $ws.Range("D16")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "53.544.46"
$rng.Style = $origStyle

$rng = $ws.Range("E16")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -5.38%  "
$rng.Style = $origStyle

$rng = $ws.Range("E17")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -4.13%  "
$rng.Style = $origStyle

$rng = $ws.Range("D18")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "2.233.66"
$rng.Style = $origStyle

$rng = $ws.Range("E18")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -6.29%  "
$rng.Style = $origStyle

$rng = $ws.Range("E19")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -4.21%  "
$rng.Style = $origStyle

$rng = $ws.Range("E20")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -2.10%  "
$rng.Style = $origStyle

$rng = $ws.Range("D21")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "296.12"
$rng.Style = $origStyle

$rng = $ws.Range("E21")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -4.40%  "
$rng.Style = $origStyle

$rng = $ws.Range("E22")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -1.89%  "
$rng.Style = $origStyle

$rng = $ws.Range("D23")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "0.997"
$rng.Style = $origStyle

$rng = $ws.Range("E23")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -0.33%  "
$rng.Style = $origStyle

$rng = $ws.Range("D24")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "63.70"
$rng.Style = $origStyle

$rng = $ws.Range("E24")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -4.02%  "
$rng.Style = $origStyle

$rng = $ws.Range("E25")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  +0.10%  "
$rng.Style = $origStyle

$rng = $ws.Range("E26")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -1.12%  "
$rng.Style = $origStyle

$rng = $ws.Range("B27")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "Kaspa"
$rng.Style = $origStyle

$rng = $ws.Range("C27")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$rng.Style = $origStyle

$rng = $ws.Range("D27")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "0.147"
$rng.Style = $origStyle

$rng = $ws.Range("E27")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -1.21%  "
$rng.Style = $origStyle

$rng = $ws.Range("B28")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "WrappedeETH"
$rng.Style = $origStyle

$rng = $ws.Range("C28")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$rng.Style = $origStyle

$rng = $ws.Range("D28")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "2.317.49"
$rng.Style = $origStyle

$rng = $ws.Range("E28")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -6.98%  "
$rng.Style = $origStyle

$rng = $ws.Range("D29")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "7.03"
$rng.Style = $origStyle

$rng = $ws.Range("E29")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -3.26%  "
$rng.Style = $origStyle

$rng = $ws.Range("D30")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "163.15"
$rng.Style = $origStyle

$rng = $ws.Range("E30")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -6.15%  "
$rng.Style = $origStyle

$rng = $ws.Range("E31")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -4.09%  "
$rng.Style = $origStyle

$rng = $ws.Range("E32")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -0.17%  "
$rng.Style = $origStyle

$rng = $ws.Range("B33")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "Aptos"
$rng.Style = $origStyle

$rng = $ws.Range("C33")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$rng.Style = $origStyle

$rng = $ws.Range("D33")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "5.79"
$rng.Style = $origStyle

$rng = $ws.Range("E33")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -1.08%  "
$rng.Style = $origStyle

$rng = $ws.Range("B34")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "PEPE"
$rng.Style = $origStyle

$rng = $ws.Range("C34")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$rng.Style = $origStyle

$rng = $ws.Range("D34")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "0.0₃0668"
$rng.Style = $origStyle

$rng = $ws.Range("E34")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -6.56%  "
$rng.Style = $origStyle

$rng = $ws.Range("D35")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "0.999"
$rng.Style = $origStyle

$rng = $ws.Range("E35")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  +0.27%  "
$rng.Style = $origStyle

$rng = $ws.Range("E36")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -2.25%  "
$rng.Style = $origStyle

$rng = $ws.Range("D37")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "17.28"
$rng.Style = $origStyle

$rng = $ws.Range("E37")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -2.30%  "
$rng.Style = $origStyle

$rng = $ws.Range("E38")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -1.14%  "
$rng.Style = $origStyle

$rng = $ws.Range("D39")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "0.833"
$rng.Style = $origStyle

$rng = $ws.Range("E39")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  +0.79%  "
$rng.Style = $origStyle

$rng = $ws.Range("E40")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -3.95%  "
$rng.Style = $origStyle

$rng = $ws.Range("E41")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -3.68%  "
$rng.Style = $origStyle

$rng = $ws.Range("E42")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -1.04%  "
$rng.Style = $origStyle

$rng = $ws.Range("E43")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -1.48%  "
$rng.Style = $origStyle

$rng = $ws.Range("D44")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "127.69"
$rng.Style = $origStyle

$rng = $ws.Range("E44")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -0.12%  "
$rng.Style = $origStyle

$rng = $ws.Range("E45")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -2.66%  "
$rng.Style = $origStyle

$rng = $ws.Range("D46")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "4.83"
$rng.Style = $origStyle

$rng = $ws.Range("E46")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  +1.41%  "
$rng.Style = $origStyle

$rng = $ws.Range("D47")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "0.0883"
$rng.Style = $origStyle

$rng = $ws.Range("E47")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -1.74%  "
$rng.Style = $origStyle

$rng = $ws.Range("D48")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "0.537"
$rng.Style = $origStyle

$rng = $ws.Range("E48")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -5.19%  "
$rng.Style = $origStyle

$rng = $ws.Range("D49")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "235.04"
$rng.Style = $origStyle

$rng = $ws.Range("E49")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -2.43%  "
$rng.Style = $origStyle

$rng = $ws.Range("E50")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -2.16%  "
$rng.Style = $origStyle

$rng = $ws.Range("E51")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = "  -3.50%  "
$rng.Style = $origStyle

Write-Output "Done applying 87 cell updates"
